# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" on the Overview sheet and the
# Priority / Correspond Handoff Datetime / Correspond Handback Datetime
# columns on the per-language report sheets (zh-cn, de-de) to reflect a
# freshly regenerated handback status report.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: refresh "Latest HO Xliff Generate Date" (column G)
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2017-01-03 04:08:43"
$overview.Range("G3").Value = "2017-01-03 04:08:43"
$overview.Range("G4").Value = "2017-01-03 04:08:43"
$overview.Range("G5").Value = "2017-01-03 04:08:43"

# ---------------------------------------------------------------------
# zh-cn sheet: Priority -> "mt", Correspond Handoff Datetime refreshed for
# all rows, Correspond Handback Datetime refreshed for rows 3 and 5
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2").Value = "mt"
$zhcn.Range("H2").Value = "2017-01-03 04:08:32"
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("H3").Value = "2017-01-03 04:08:32"
$zhcn.Range("L3").Value = "2017-01-03 04:09:05"
$zhcn.Range("E4").Value = "mt"
$zhcn.Range("H4").Value = "2017-01-03 04:08:32"
$zhcn.Range("E5").Value = "mt"
$zhcn.Range("H5").Value = "2017-01-03 04:08:32"
$zhcn.Range("L5").Value = "2017-01-03 04:09:05"

# ---------------------------------------------------------------------
# de-de sheet: Priority -> "mt", Correspond Handoff Datetime refreshed for
# all rows, Correspond Handback Datetime refreshed for rows 3 and 5
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2").Value = "mt"
$dede.Range("H2").Value = "2017-01-03 04:08:43"
$dede.Range("E3").Value = "mt"
$dede.Range("H3").Value = "2017-01-03 04:08:43"
$dede.Range("L3").Value = "2017-01-03 04:09:16"
$dede.Range("E4").Value = "mt"
$dede.Range("H4").Value = "2017-01-03 04:08:43"
$dede.Range("E5").Value = "mt"
$dede.Range("H5").Value = "2017-01-03 04:08:43"
$dede.Range("L5").Value = "2017-01-03 04:09:16"
